# Natmi following Dr Hou advice
# Adds an "M2" sending/target cluster category to the Gal-Gpr151 ligand-receptor
# table and refreshes all of the computed expression/specificity statistics.
# The data block (rows 2-13) is fully replaced: previously there were only two
# cluster categories (FAPs, sCs); now there are three (FAPs, M2, sCs), combined
# with each other and laid out in A..T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("FAPs", "Gal", "Gpr151", "ECs", 3, 1, 1.642758, 4.928274, 0.08311198773343359, 0.08311198773343359, 3, 1, 0.271848, 0.815544, 0.1256018707893687, 0.1256018707893686, 0.446580476784, 4.019224291056, 0.01043902114434232, 0.01043902114434232),
  @("FAPs", "Gal", "Gpr151", "FAPs", 3, 1, 1.642758, 4.928274, 0.08311198773343359, 0.08311198773343359, 3, 1, 0.5236423333333333, 1.570927, 0.2419383504428094, 0.2419383504428093, 0.860217632222, 7.741958689998, 0.02010797721424993, 0.02010797721424993),
  @("FAPs", "Gal", "Gpr151", "sCs", 3, 1, 1.642758, 4.928274, 0.08311198773343359, 0.08311198773343359, 3, 1, 0.323177, 0.9695309999999999, 0.1493173972076152, 0.1493173972076152, 0.530901602166, 4.778114419494, 0.01241006568510755, 0.01241006568510755),
  @("FAPs", "Gal", "Gpr151", "M2", 3, 1, 1.642758, 4.928274, 0.08311198773343359, 0.08311198773343359, 3, 1, 1.045695333333333, 3.137086, 0.4831423815602068, 0.4831423815602067, 1.717824374396, 15.460419369564, 0.0401549236897338, 0.04015492368973379),
  @("M2", "Gal", "Gpr151", "ECs", 2, 0.6666666666666666, 0.4045783333333333, 1.213735, 0.02046881493026545, 0.02046881493026545, 3, 1, 0.271848, 0.815544, 0.1256018707893687, 0.1256018707893686, 0.10998381076, 0.98985429684, 0.002570921448082701, 0.0025709214480827),
  @("M2", "Gal", "Gpr151", "FAPs", 2, 0.6666666666666666, 0.4045783333333333, 1.213735, 0.02046881493026545, 0.02046881493026545, 3, 1, 0.5236423333333333, 1.570927, 0.2419383504428094, 0.2419383504428093, 0.2118543424827778, 1.906689082345, 0.00495219131974757, 0.004952191319747569),
  @("M2", "Gal", "Gpr151", "sCs", 2, 0.6666666666666666, 0.4045783333333333, 1.213735, 0.02046881493026545, 0.02046881493026545, 3, 1, 0.323177, 0.9695309999999999, 0.1493173972076152, 0.1493173972076152, 0.1307504120316667, 1.176753708285, 0.003056350169311611, 0.003056350169311611),
  @("M2", "Gal", "Gpr151", "M2", 2, 0.6666666666666666, 0.4045783333333333, 1.213735, 0.02046881493026545, 0.02046881493026545, 3, 1, 1.045695333333333, 3.137086, 0.4831423815602068, 0.4831423815602067, 0.4230656751344444, 3.80759107621, 0.009889351993123565, 0.009889351993123564),
  @("sCs", "Gal", "Gpr151", "ECs", 3, 1, 17.71826, 53.15478, 0.8964191973363009, 0.8964191973363009, 3, 1, 0.271848, 0.815544, 0.1256018707893687, 0.1256018707893686, 4.81667354448, 43.35006190032001, 0.1125919281969436, 0.1125919281969436),
  @("sCs", "Gal", "Gpr151", "FAPs", 3, 1, 17.71826, 53.15478, 0.8964191973363009, 0.8964191973363009, 3, 1, 0.5236423333333333, 1.570927, 0.2419383504428094, 0.2419383504428093, 9.278031009006666, 83.50227908106, 0.2168781819088119, 0.2168781819088118),
  @("sCs", "Gal", "Gpr151", "sCs", 3, 1, 17.71826, 53.15478, 0.8964191973363009, 0.8964191973363009, 3, 1, 0.323177, 0.9695309999999999, 0.1493173972076152, 0.1493173972076152, 5.72613411202, 51.53520700818, 0.1338509813531961, 0.1338509813531961),
  @("sCs", "Gal", "Gpr151", "M2", 3, 1, 17.71826, 53.15478, 0.8964191973363009, 0.8964191973363009, 3, 1, 1.045695333333333, 3.137086, 0.4831423815602068, 0.4831423815602067, 18.52790179678667, 166.75111617108, 0.4330981058773494, 0.4330981058773494)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item(2 + $r, 1 + $c).Value = $row[$c]
    }
}
